# Implemented Custom Excel Report
# Adds a new "Negative" worksheet at the end of the workbook containing an
# "Email" header and a mailto hyperlink to a test email address, and makes
# it the active/selected sheet (mirroring the eCare sheet losing focus).

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Negative"

# Populate the content
$newSheet.Range("A1").Value = "Email"
$newSheet.Range("A2").Value = "ScriptFaile@gmail.com"

# Turn the email address into a mailto hyperlink (this also applies the
# built-in "Hyperlink" cell style / underlined themed font automatically)
[void]$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "mailto:ScriptFaile@gmail.com")

# Make the new sheet the active one with A2 selected
[void]$newSheet.Range("A2").Select()
[void]$newSheet.Activate()
